# Insert a new data row at row 268 (Hortaliza, Femacal de La Calera - Cilantro),
# shifting existing rows 268:322 down to 269:323, then populate the new row
# with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 268..322 down by one row.
$ws.Rows.Item(268).Insert()

# Populate the newly-inserted row 268 with its data.
$ws.Cells.Item(268, 1).Value = 3
$ws.Cells.Item(268, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(268, 3).Value = "Coquimbo"
$ws.Cells.Item(268, 4).Value = 44637
$ws.Cells.Item(268, 5).Value = 5
$ws.Cells.Item(268, 6).Value = 100112040
$ws.Cells.Item(268, 7).Value = "Cilantro"
$ws.Cells.Item(268, 8).Value = "Sin especificar"
$ws.Cells.Item(268, 9).Value = "Primera"
$ws.Cells.Item(268, 10).Value = 120
$ws.Cells.Item(268, 11).Value = 4500
$ws.Cells.Item(268, 12).Value = 5000
$ws.Cells.Item(268, 13).Value = 4750
$ws.Cells.Item(268, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(268, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(268, 16).Value = 1583
$ws.Cells.Item(268, 17).Value = 3
$ws.Cells.Item(268, 18).Value = "Hortaliza"
